$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 722.44446
$ws.Range("I9").Value = 1200
$ws.Range("J9").Value = 125.5
$ws.Range("K9").Value = 1200
$ws.Range("L9").Value = 125.5
$ws.Range("M9").Value = -1031
$ws.Range("N9").Value = -463.5
$ws.Range("H17").Value = 5188.4614
$ws.Range("J17").Value = 6223.8096
$ws.Range("L17").Value = 18671.4288
$ws.Range("N17").Value = -19007.4288
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = 0
$ws.Range("H38").Value = 516
$ws.Range("I38").Value = 141.875
$ws.Range("J38").Value = 2012.5
$ws.Range("K38").Value = 425.625
$ws.Range("L38").Value = 6037.5
$ws.Range("M38").Value = -53.625
$ws.Range("N38").Value = -6781.5
$ws.Range("H97").Value = 956.6667
$ws.Range("J97").Value = 2000
$ws.Range("L97").Value = 6000
$ws.Range("N97").Value = -6992
$ws.Range("H131").Value = 7219.6
$ws.Range("I131").Value = 5897.3335
$ws.Range("K131").Value = 17692.0005
$ws.Range("M131").Value = -12652.0005
$ws.Range("H135").Value = 743
$ws.Range("I135").Value = 597.4583
$ws.Range("J135").Value = 1325.1666
$ws.Range("K135").Value = 5377.1247
$ws.Range("L135").Value = 11926.4994
$ws.Range("M135").Value = -2842.1247
$ws.Range("N135").Value = -16996.4994
$ws.Range("H137").Value = 50004108
$ws.Range("I137").Value = 200004940
$ws.Range("J137").Value = 3829.0667
$ws.Range("K137").Value = 600014820
$ws.Range("L137").Value = 11487.2001
$ws.Range("M137").Value = -600012270
$ws.Range("N137").Value = -16587.2001

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 11495482
$ws.Range("I74").Value = 11495482
$ws.Range("K74").Value = 11495482
$ws.Range("M74").Value = -11494608
$ws.Range("H77").Value = 11495482
$ws.Range("I77").Value = 11495482
$ws.Range("K77").Value = 57477410
$ws.Range("M77").Value = -57473042
$ws.Range("H122").Value = 125001390
$ws.Range("I122").Value = 1491
$ws.Range("K122").Value = 4473
$ws.Range("M122").Value = -2023
$ws.Range("H132").Value = 3106.3809
$ws.Range("I132").Value = 2074.5
$ws.Range("J132").Value = 12909.25
$ws.Range("K132").Value = 6223.5
$ws.Range("L132").Value = 38727.75
$ws.Range("M132").Value = -3693.5
$ws.Range("N132").Value = -43787.75

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1465
$ws.Range("I99").Value = 1465
$ws.Range("K99").Value = 1465
$ws.Range("M99").Value = 33
$ws.Range("H134").Value = 4112.4375
$ws.Range("I134").Value = 1752.6923
$ws.Range("J134").Value = 14338
$ws.Range("K134").Value = 5258.0769
$ws.Range("L134").Value = 43014
$ws.Range("M134").Value = -2723.0769
$ws.Range("N134").Value = -48084

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 178671.33
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 178671.33
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 178671.33
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -179261.33
$ws.Range("H34").Value = 178671.33
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 178671.33
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 178671.33
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -179075.33
$ws.Range("H58").Value = 4553.9546
$ws.Range("I58").Value = 1960.0834
$ws.Range("J58").Value = 7666.6
$ws.Range("K58").Value = 1960.0834
$ws.Range("L58").Value = 7666.6
$ws.Range("M58").Value = -1757.0834
$ws.Range("N58").Value = -8072.6
$ws.Range("H107").Value = 1030.2778
$ws.Range("I107").Value = 942.9167
$ws.Range("J107").Value = 1205
$ws.Range("K107").Value = 942.9167
$ws.Range("L107").Value = 1205
$ws.Range("M107").Value = 977.0833
$ws.Range("N107").Value = -5045
$ws.Range("H122").Value = 6368.55
$ws.Range("I122").Value = 2836.2144
$ws.Range("K122").Value = 8508.643199999999
$ws.Range("M122").Value = -6058.643199999999
$ws.Range("H132").Value = 1871.1428
$ws.Range("I132").Value = 552.46155
$ws.Range("J132").Value = 19014
$ws.Range("K132").Value = 1657.38465
$ws.Range("L132").Value = 57042
$ws.Range("M132").Value = 872.61535
$ws.Range("N132").Value = -62102
$ws.Range("H136").Value = 4553.9546
$ws.Range("I136").Value = 1960.0834
$ws.Range("J136").Value = 7666.6
$ws.Range("K136").Value = 5880.2502
$ws.Range("L136").Value = 22999.8
$ws.Range("M136").Value = -3330.2502
$ws.Range("N136").Value = -28099.8

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11993698
$ws.Range("I4").Value = 5714384
$ws.Range("K4").Value = 17143152
$ws.Range("M4").Value = -17143040
$ws.Range("H5").Value = 10332.375
$ws.Range("J5").Value = 19751
$ws.Range("L5").Value = 59253
$ws.Range("N5").Value = -59477
$ws.Range("H23").Value = 242.63637
$ws.Range("I23").Value = 196.4
$ws.Range("J23").Value = 281.16666
$ws.Range("K23").Value = 589.2
$ws.Range("L23").Value = 843.4999799999999
$ws.Range("M23").Value = -354.2
$ws.Range("N23").Value = -1313.49998
$ws.Range("H132").Value = 4068.182
$ws.Range("I132").Value = 2458.4167
$ws.Range("J132").Value = 5999.9
$ws.Range("K132").Value = 22125.7503
$ws.Range("L132").Value = 53999.1
$ws.Range("M132").Value = -19595.7503
$ws.Range("N132").Value = -59059.1
$ws.Range("H135").Value = 10332.375
$ws.Range("J135").Value = 19751
$ws.Range("L135").Value = 177759
$ws.Range("N135").Value = -182829

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 24723.75
$ws.Range("J33").Value = 24723.75
$ws.Range("L33").Value = 24723.75
$ws.Range("N33").Value = -25227.75
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H132").Value = 4373444.5
$ws.Range("I132").Value = 4373444.5
$ws.Range("K132").Value = 13120333.5
$ws.Range("M132").Value = -13117803.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 2449.75
$ws.Range("J35").Value = 2713.4
$ws.Range("L35").Value = 2713.4
$ws.Range("N35").Value = -3385.4
$ws.Range("H122").Value = 4521.1113
$ws.Range("I122").Value = 1377
$ws.Range("J122").Value = 8451.25
$ws.Range("K122").Value = 4131
$ws.Range("L122").Value = 25353.75
$ws.Range("M122").Value = -1681
$ws.Range("N122").Value = -30253.75
$ws.Range("H132").Value = 3762.3389
$ws.Range("I132").Value = 1986.3334
$ws.Range("J132").Value = 8150.1177
$ws.Range("K132").Value = 5959.0002
$ws.Range("L132").Value = 24450.3531
$ws.Range("M132").Value = -3429.0002
$ws.Range("N132").Value = -29510.3531
$ws.Range("H136").Value = 7770.227
$ws.Range("I136").Value = 2929.4167
$ws.Range("K136").Value = 8788.250100000001
$ws.Range("M136").Value = -6238.250100000001

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H124").Value = 69161
$ws.Range("I124").Value = 65390
$ws.Range("J124").Value = 70418
$ws.Range("K124").Value = 65390
$ws.Range("L124").Value = 70418
$ws.Range("M124").Value = -60480
$ws.Range("N124").Value = -80238
$ws.Range("H132").Value = 7218.4897
$ws.Range("I132").Value = 5336.0513
$ws.Range("J132").Value = 14560
$ws.Range("K132").Value = 16008.1539
$ws.Range("L132").Value = 43680
$ws.Range("M132").Value = -13478.1539
$ws.Range("N132").Value = -48740
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

Write-Host "Applied all cell updates."